$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 8065498.5
$ws.Range("I135").Value = 652.0465
$ws.Range("J135").Value = 26317520
$ws.Range("K135").Value = 5868.418500000001
$ws.Range("L135").Value = 236857680
$ws.Range("M135").Value = -3333.418500000001
$ws.Range("N135").Value = -236862750

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2543.92
$ws.Range("I138").Value = 976.62164
$ws.Range("J138").Value = 3464.3967
$ws.Range("K138").Value = 2929.86492
$ws.Range("L138").Value = 10393.1901
$ws.Range("M138").Value = 2210.13508
$ws.Range("N138").Value = -20673.1901

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1933.28
$ws.Range("I32").Value = 1326.069
$ws.Range("J32").Value = 5996.923
$ws.Range("K32").Value = 1326.069
$ws.Range("L32").Value = 5996.923
$ws.Range("M32").Value = -1039.069
$ws.Range("N32").Value = -6570.923

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1329.0132
$ws.Range("I61").Value = 1197.2549
$ws.Range("J61").Value = 1597.8
$ws.Range("K61").Value = 1197.2549
$ws.Range("L61").Value = 1597.8
$ws.Range("M61").Value = -985.2548999999999
$ws.Range("N61").Value = -2021.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1779.221
$ws.Range("I74").Value = 1648.4084
$ws.Range("J74").Value = 2398.4
$ws.Range("K74").Value = 1648.4084
$ws.Range("L74").Value = 2398.4
$ws.Range("M74").Value = -774.4084
$ws.Range("N74").Value = -4146.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1779.221
$ws.Range("I77").Value = 1648.4084
$ws.Range("J77").Value = 2398.4
$ws.Range("K77").Value = 8242.041999999999
$ws.Range("L77").Value = 11992
$ws.Range("M77").Value = -3874.041999999999
$ws.Range("N77").Value = -20728

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H109").Value = 47894.832
$ws.Range("J109").Value = 47894.832
$ws.Range("L109").Value = 47894.832
$ws.Range("N109").Value = -50668.832

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H117").Value = 36938.715
$ws.Range("J117").Value = 36938.715
$ws.Range("L117").Value = 36938.715
$ws.Range("N117").Value = -46116.715

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 7693846
$ws.Range("I132").Value = 12196138
$ws.Range("J132").Value = 2430.7917
$ws.Range("K132").Value = 36588414
$ws.Range("L132").Value = 7292.375100000001
$ws.Range("M132").Value = -36585884
$ws.Range("N132").Value = -12352.3751

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1329.0132
$ws.Range("I136").Value = 1197.2549
$ws.Range("J136").Value = 1597.8
$ws.Range("K136").Value = 3591.7647
$ws.Range("L136").Value = 4793.4
$ws.Range("M136").Value = -1041.7647
$ws.Range("N136").Value = -9893.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H133").Value = 48747.75
$ws.Range("J133").Value = 48747.75
$ws.Range("L133").Value = 48747.75
$ws.Range("N133").Value = -58867.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2272.89
$ws.Range("I134").Value = 903.8936
$ws.Range("J134").Value = 3486.9058
$ws.Range("K134").Value = 2711.6808
$ws.Range("L134").Value = 10460.7174
$ws.Range("M134").Value = -176.6808000000001
$ws.Range("N134").Value = -15530.7174

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2891.71
$ws.Range("I31").Value = 1072
$ws.Range("J31").Value = 3564.7534
$ws.Range("K31").Value = 1072
$ws.Range("L31").Value = 3564.7534
$ws.Range("M31").Value = -777
$ws.Range("N31").Value = -4154.7534

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2891.71
$ws.Range("I34").Value = 1072
$ws.Range("J34").Value = 3564.7534
$ws.Range("K34").Value = 1072
$ws.Range("L34").Value = 3564.7534
$ws.Range("M34").Value = -870
$ws.Range("N34").Value = -3968.7534

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1524.1082
$ws.Range("I58").Value = 1128
$ws.Range("J58").Value = 2174.8572
$ws.Range("K58").Value = 1128
$ws.Range("L58").Value = 2174.8572
$ws.Range("M58").Value = -925
$ws.Range("N58").Value = -2580.8572

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H118").Value = 48724.668
$ws.Range("J118").Value = 48724.668
$ws.Range("L118").Value = 48724.668
$ws.Range("N118").Value = -52038.668

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 32415.11
$ws.Range("I132").Value = 1488.3784
$ws.Range("J132").Value = 159558.33
$ws.Range("K132").Value = 4465.135200000001
$ws.Range("L132").Value = 478674.99
$ws.Range("M132").Value = -1935.135200000001
$ws.Range("N132").Value = -483734.99

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 287242.22
$ws.Range("I134").Value = 974.34283
$ws.Range("J134").Value = 1002911.94
$ws.Range("K134").Value = 2923.02849
$ws.Range("L134").Value = 3008735.82
$ws.Range("M134").Value = -388.0284900000001
$ws.Range("N134").Value = -3013805.82

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1524.1082
$ws.Range("I136").Value = 1128
$ws.Range("J136").Value = 2174.8572
$ws.Range("K136").Value = 3384
$ws.Range("L136").Value = 6524.571599999999
$ws.Range("M136").Value = -834
$ws.Range("N136").Value = -11624.5716

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 4144.6787
$ws.Range("I113").Value = 6404.1763
$ws.Range("J113").Value = 652.7273
$ws.Range("K113").Value = 19212.5289
$ws.Range("L113").Value = 1958.1819
$ws.Range("M113").Value = -17042.5289
$ws.Range("N113").Value = -6298.1819

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H110").Value = 47662
$ws.Range("J110").Value = 47662
$ws.Range("L110").Value = 47662
$ws.Range("N110").Value = -55842

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1369.7894
$ws.Range("I113").Value = 1336.4546
$ws.Range("K113").Value = 1336.4546
$ws.Range("M113").Value = 833.5454

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1555.4445
$ws.Range("I122").Value = 1599.875
$ws.Range("J122").Value = 1200
$ws.Range("K122").Value = 4799.625
$ws.Range("L122").Value = 3600
$ws.Range("M122").Value = -2349.625
$ws.Range("N122").Value = -8500

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 6339.3335
$ws.Range("I126").Value = 9789.666999999999
$ws.Range("J126").Value = 1738.8889
$ws.Range("K126").Value = 29369.001
$ws.Range("L126").Value = 5216.6667
$ws.Range("M126").Value = -26899.001
$ws.Range("N126").Value = -10156.6667

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2791.7742
$ws.Range("I132").Value = 2042.5
$ws.Range("J132").Value = 3408.8235
$ws.Range("K132").Value = 6127.5
$ws.Range("L132").Value = 10226.4705
$ws.Range("M132").Value = -3597.5
$ws.Range("N132").Value = -15286.4705

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 734.2632
$ws.Range("I22").Value = 731.5833
$ws.Range("J22").Value = 738.8570999999999
$ws.Range("K22").Value = 731.5833
$ws.Range("L22").Value = 738.8570999999999
$ws.Range("M22").Value = -436.5833
$ws.Range("N22").Value = -1328.8571

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 734.2632
$ws.Range("I27").Value = 731.5833
$ws.Range("J27").Value = 738.8570999999999
$ws.Range("K27").Value = 731.5833
$ws.Range("L27").Value = 738.8570999999999
$ws.Range("M27").Value = -624.5833
$ws.Range("N27").Value = -952.8570999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 8075.353
$ws.Range("I46").Value = 3662.375
$ws.Range("J46").Value = 11998
$ws.Range("K46").Value = 3662.375
$ws.Range("L46").Value = 11998
$ws.Range("M46").Value = -3474.375
$ws.Range("N46").Value = -12374

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3187.5
$ws.Range("I68").Value = 3000
$ws.Range("J68").Value = 3250
$ws.Range("K68").Value = 3000
$ws.Range("L68").Value = 3250
$ws.Range("M68").Value = -2251
$ws.Range("N68").Value = -4748

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 3187.5
$ws.Range("I71").Value = 3000
$ws.Range("J71").Value = 3250
$ws.Range("K71").Value = 15000
$ws.Range("L71").Value = 16250
$ws.Range("M71").Value = -11256
$ws.Range("N71").Value = -23738

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 60901.234
$ws.Range("I122").Value = 101771.6
$ws.Range("J122").Value = 2515
$ws.Range("K122").Value = 305314.8
$ws.Range("L122").Value = 7545
$ws.Range("M122").Value = -302864.8
$ws.Range("N122").Value = -12445

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2703.5642
$ws.Range("I132").Value = 1810.3478
$ws.Range("K132").Value = 5431.0434
$ws.Range("M132").Value = -2901.0434

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 1622.7805
$ws.Range("I136").Value = 1335.3429
$ws.Range("J136").Value = 3299.5
$ws.Range("K136").Value = 4006.0287
$ws.Range("L136").Value = 9898.5
$ws.Range("M136").Value = -1456.0287
$ws.Range("N136").Value = -14998.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 398.86667
$ws.Range("I100").Value = 365.25
$ws.Range("K100").Value = 730.5
$ws.Range("M100").Value = -189.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 7143557
$ws.Range("I122").Value = 9524409
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 28573227
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = -28570777
$ws.Range("N122").Value = -7900

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2057.9075
$ws.Range("I132").Value = 2011.3684
$ws.Range("J132").Value = 2168.4375
$ws.Range("K132").Value = 6034.1052
$ws.Range("L132").Value = 6505.3125
$ws.Range("M132").Value = -3504.1052
$ws.Range("N132").Value = -11565.3125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H133").Value = 64233.25
$ws.Range("J133").Value = 64233.25
$ws.Range("L133").Value = 64233.25
$ws.Range("N133").Value = -74353.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 927.3125
$ws.Range("I136").Value = 698.0323
$ws.Range("J136").Value = 1142.697
$ws.Range("K136").Value = 2094.0969
$ws.Range("L136").Value = 3428.090999999999
$ws.Range("M136").Value = 455.9031
$ws.Range("N136").Value = -8528.091

